# Applies highlighting edits to CapstoneUserStories.docx
#
# WdColorIndex values observed from this runtime's <w:highlight> mapping:
#   4 = green (wdBrightGreen)  -> w:val="green"
#   6 = red   (wdRed)          -> w:val="red"
#   7 = yellow(wdYellow)       -> w:val="yellow"
#
# NOTE: Setting HighlightColorIndex on a Range's *.Font* sub-object (rather
# than directly on the Range) is what makes the paragraph-mark's rPr (inside
# w:pPr) pick up the highlight too, matching Word's real behavior of
# stamping the mark when the whole paragraph (incl. its end) is selected.

$d = $word.ActiveDocument

function Highlight($startPos, $endPos, $colorIndex) {
    $rng = $d.Range($startPos, $endPos)
    $rng.Font.HighlightColorIndex = $colorIndex
}

function HighlightWholeParagraph($paraIndex, $colorIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Font.HighlightColorIndex = $colorIndex
}

# --- Paragraph 3: "Entity Relationship Diagram ... 25 weighted points  " ---
# Whole paragraph (incl. mark) -> green, matching the trailing "  " run which
# then gets re-cleared to no-highlight by re-setting just that tail range.
HighlightWholeParagraph 3 4
Highlight 162 165 0

# --- Paragraph 7: "(5 points) As a Job Seeker ... login and logout." ---
# Paragraph mark -> yellow; most of the text -> green; final " logout." -> yellow.
HighlightWholeParagraph 7 7
Highlight 228 337 4

# --- Paragraph 8: "(5 points) As a Job Seeker ... zip code." ---
# Whole paragraph (incl. mark) -> red.
HighlightWholeParagraph 8 6

# --- Paragraph 9: "(5 points) As a Job Seeker ... apply to the jobs I choose." ---
# Whole paragraph (incl. mark) -> red.
HighlightWholeParagraph 9 6

# --- Paragraph 10: "(2.5 points) As a Job Poster ... create job postings" ---
# Whole paragraph (incl. mark) -> yellow.
HighlightWholeParagraph 10 7

# --- Paragraph 11: "(7.5 points) As a Job Poster ... for my job." ---
# Whole paragraph (incl. mark) -> red.
HighlightWholeParagraph 11 6

# --- Paragraph 12: "(2.5 points) As a Job Seeker ... applied to." ---
# Whole paragraph (incl. mark) -> red.
HighlightWholeParagraph 12 6

# --- Paragraph 13: "(10 points) As a Job Seeker ... Directions API)." ---
# Whole paragraph (incl. mark) -> red.
HighlightWholeParagraph 13 6

# --- Paragraph 14: "(10 points) As a developer ... etc.)" ---
# Whole paragraph (incl. mark) -> red; text also gets split into three runs
# around "experience" (same red highlight on all three, so the split is a
# side-effect we create explicitly to mirror the target run layout).
HighlightWholeParagraph 14 6
Highlight 1021 1031 6

Write-Host "done"
